$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "合富中国"
$ws.Range("B2").Value = "海马汽车"
$ws.Range("C2").Value = "平潭发展"

$ws.Range("A3").Value = "中国中免"
$ws.Range("B3").Value = "中国中免"
$ws.Range("C3").Value = "海马汽车"

$ws.Range("A4").Value = "海马汽车"
$ws.Range("B4").Value = "合富中国"
$ws.Range("C4").Value = "合富中国"

$ws.Range("A5").Value = "平潭发展"
$ws.Range("B5").Value = "平潭发展"
$ws.Range("C5").Value = "香农芯创"

$ws.Range("A6").Value = "多氟多"
$ws.Range("B6").Value = "多氟多"
$ws.Range("C6").Value = "雪人集团"

$ws.Range("A7").Value = "天际股份"
$ws.Range("B7").Value = "特变电工"
$ws.Range("C7").Value = "中国中免"

$ws.Range("A8").Value = "四川金顶"
$ws.Range("B8").Value = "天际股份"
$ws.Range("C8").Value = "天际股份"

$ws.Range("A9").Value = "闻泰科技"
$ws.Range("B9").Value = "闻泰科技"
$ws.Range("C9").Value = "盈新发展"

$ws.Range("A10").Value = "雪人集团"
$ws.Range("B10").Value = "雪人集团"
$ws.Range("C10").Value = "特变电工"

$ws.Range("A11").Value = "特变电工"
$ws.Range("B11").Value = "四川金顶"
$ws.Range("C11").Value = "多氟多"

$ws.Range("A12").Value = "香农芯创"
$ws.Range("B12").Value = "香农芯创"
$ws.Range("C12").Value = "摩恩电气"

$ws.Range("A13").Value = "盈新发展"
$ws.Range("B13").Value = "中国西电"
$ws.Range("C13").Value = "海陆重工"

$ws.Range("A14").Value = "三花智控"
$ws.Range("B14").Value = "三花智控"
$ws.Range("C14").Value = "闻泰科技"

$ws.Range("A15").Value = "海陆重工"
$ws.Range("B15").Value = "盈新发展"
$ws.Range("C15").Value = "三花智控"

$ws.Range("A16").Value = "保变电气"
$ws.Range("B16").Value = "保变电气"
$ws.Range("C16").Value = "四川金顶"

$ws.Range("A17").Value = "摩恩电气"
$ws.Range("B17").Value = "海陆重工"
$ws.Range("C17").Value = "清水源"

$ws.Range("A18").Value = "中国西电"
$ws.Range("B18").Value = "摩恩电气"
$ws.Range("C18").Value = "厦工股份"

$ws.Range("A19").Value = "工业富联"
$ws.Range("B19").Value = "工业富联"
$ws.Range("C19").Value = "保变电气"

$ws.Range("A20").Value = "普天科技"
$ws.Range("B20").Value = "澄星股份"
$ws.Range("C20").Value = "合肥城建"

$ws.Range("A21").Value = "乾照光电"
$ws.Range("B21").Value = "天赐材料"
$ws.Range("C21").Value = "工业富联"
